$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 29-30 (shifts old rows 29-47 down to 31-49)
$ws.Range("A29:A30").EntireRow.Insert()

# Row 29: new Tuna / Primera record
$ws.Cells.Item(29, 1).Value = 1
$ws.Cells.Item(29, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(29, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(29, 4).Value = 44574
$ws.Cells.Item(29, 5).Value = 15
$ws.Cells.Item(29, 6).Value = 100112027
$ws.Cells.Item(29, 7).Value = "Melón"
$ws.Cells.Item(29, 8).Value = "Tuna"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 50
$ws.Cells.Item(29, 11).Value = 6500
$ws.Cells.Item(29, 12).Value = 7000
$ws.Cells.Item(29, 13).Value = 6750
$ws.Cells.Item(29, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(29, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(29, 16).Value = 375
$ws.Cells.Item(29, 17).Value = 18
$ws.Cells.Item(29, 18).Value = "Hortaliza"

# Row 30: new Tuna / Segunda record
$ws.Cells.Item(30, 1).Value = 1
$ws.Cells.Item(30, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(30, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(30, 4).Value = 44574
$ws.Cells.Item(30, 5).Value = 15
$ws.Cells.Item(30, 6).Value = 100112027
$ws.Cells.Item(30, 7).Value = "Melón"
$ws.Cells.Item(30, 8).Value = "Tuna"
$ws.Cells.Item(30, 9).Value = "Segunda"
$ws.Cells.Item(30, 10).Value = 70
$ws.Cells.Item(30, 11).Value = 4500
$ws.Cells.Item(30, 12).Value = 5000
$ws.Cells.Item(30, 13).Value = 4750
$ws.Cells.Item(30, 14).Value = "$/caja 24 unidades"
$ws.Cells.Item(30, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(30, 16).Value = 198
$ws.Cells.Item(30, 17).Value = 24
$ws.Cells.Item(30, 18).Value = "Hortaliza"

# Make sure the new date cells keep the same date number format as column D elsewhere
$ws.Cells.Item(29, 4).NumberFormat = $ws.Cells.Item(28, 4).NumberFormat
$ws.Cells.Item(30, 4).NumberFormat = $ws.Cells.Item(28, 4).NumberFormat
